$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.458.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.104.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5216"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.56"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08913"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.114.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.812"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.992"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06649"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.312"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.516.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.347"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.349.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.531"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.204"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.391"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.946"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.802"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02579"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06831"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2302"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6868"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.320"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.664"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  +21.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.248"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
